$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.875.39'
$ws.Range("E2").Value = '  -1.36%  '

$ws.Range("D3").Value = '1.994.68'
$ws.Range("E3").Value = '  -2.93%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.74%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.605'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.73%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.16'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.37%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.376'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0780'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.104'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.92%  '

$ws.Range("D12").Value = '2.289.60'
$ws.Range("E12").Value = '  -2.81%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.03'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.12%  '

$ws.Range("E14").Value = '  -3.73%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.735'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.08'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.14%  '

$ws.Range("D17").Value = '1.999.37'
$ws.Range("E17").Value = '  -2.80%  '

$ws.Range("D18").Value = '36.769.01'
$ws.Range("E18").Value = '  -1.42%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.57%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.62%  '

$ws.Range("D21").Value = '0.0₃0811'
$ws.Range("E21").Value = '  -1.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '222.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.96%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.05%  '

$ws.Range("E24").Value = '  +0.75%  '

$ws.Range("E25").Value = '  -7.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.85'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.96%  '

$ws.Range("E28").Value = '  -4.63%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.49%  '

$ws.Range("E30").Value = '  -0.77%  '

$ws.Range("E31").Value = '  -3.85%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.47'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.45%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0606'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.73%  '

$ws.Range("E34").Value = '  -4.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.32'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.89%  '

$ws.Range("E36").Value = '  +1.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.12'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.75%  '

$ws.Range("E39").Value = '  -0.80%  '

$ws.Range("D40").Value = '1.463.26'
$ws.Range("E40").Value = '  -1.08%  '

$ws.Range("E41").Value = '  -4.60%  '

$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '94.23'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.23%  '

$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.33'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.87%  '

$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0913'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.44%  '

$ws.Range("E45").Value = '  -4.87%  '

$ws.Range("E46").Value = '  -6.47%  '

$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.99%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.11'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.90'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.72%  '

$ws.Range("D50").Value = '2.180.25'
$ws.Range("E50").Value = '  -2.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.14'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.68%  '
